$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant_vocab_mapping")

# Quick change for taxonomy value: B17 goes from numeric 1566 to text "ENERGY DATA.INFO"
$ws.Range("B17").Value = "ENERGY DATA.INFO"
